$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.797.05"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "2.082.58"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").Value = "  +0.05%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "53.87"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -6.76%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "58.86"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("E10").Value = "  -3.89%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0762"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("E12").Value = "  +1.29%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "14.95"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -5.73%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.884"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +2.84%  "
$ws.Range("D15").Value = "2.382.91"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("E16").Value = "  -3.61%  "
$ws.Range("D17").Value = "2.058.34"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "36.769.40"
$ws.Range("E18").Value = "  -0.95%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "17.21"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -3.43%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "72.66"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("D21").Value = "0.0₃0878"
$ws.Range("E21").Value = "  -1.32%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.45"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +2.03%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "238.41"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("E24").Value = "  +0.10%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.40"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("E26").Value = "  +3.55%  "
$ws.Range("E27").Value = "  -0.65%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "167.41"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.52%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "20.58"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("E30").Value = "  -0.92%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "5.31"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +10.87%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.18"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +6.38%  "
$ws.Range("E33").Value = "  +4.41%  "
$ws.Range("E34").Value = "  -1.45%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.38"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +6.47%  "
$ws.Range("E36").Value = "  +0.20%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.85"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +4.47%  "
$ws.Range("E38").Value = "  -6.79%  "
$ws.Range("E39").Value = "  -4.36%  "
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("E42").Value = "  -6.19%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.0948"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.86%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "96.36"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.76%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.85"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -12.83%  "
$ws.Range("E46").Value = "  -6.71%  "
$ws.Range("D47").Value = "1.352.40"
$ws.Range("E47").Value = "  +6.15%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.44"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.22%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "7.24"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +6.61%  "
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("D51").Value = "2.265.94"
$ws.Range("E51").Value = "  +1.53%  "
